$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value would otherwise be auto-parsed as a
# number by Excel (losing significant trailing zeros / exact text form).
# Force them to Text format first so the literal string is preserved,
# matching the original file which stores these as text.
$textFormatCells = @(
    "D5",
    "D8",
    "D9",
    "D10",
    "D14",
    "D15",
    "D18",
    "D20",
    "D21",
    "D22",
    "D23",
    "D24",
    "D25",
    "D26",
    "D28",
    "D29",
    "D31",
    "D32",
    "D34",
    "D37",
    "D39",
    "D40",
    "D41",
    "D44",
    "D46",
)
foreach ($cellRef in $textFormatCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.015.01"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "1.630.72"
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("E4").Value = "  +0.48%  "
$ws.Range("D5").Value = "214.31"
$ws.Range("E5").Value = "  -0.96%  "
$ws.Range("E6").Value = "  -1.11%  "
$ws.Range("E7").Value = "  +0.47%  "
$ws.Range("D8").Value = "0.249"
$ws.Range("E8").Value = "  -2.81%  "
$ws.Range("D9").Value = "0.0618"
$ws.Range("E9").Value = "  -3.42%  "
$ws.Range("D10").Value = "18.30"
$ws.Range("E10").Value = "  -6.88%  "
$ws.Range("E11").Value = "  -0.72%  "
$ws.Range("D12").Value = "1.859.09"
$ws.Range("E12").Value = "  -0.69%  "
$ws.Range("D13").Value = "1.622.56"
$ws.Range("E13").Value = "  -3.13%  "
$ws.Range("D14").Value = "4.16"
$ws.Range("E14").Value = "  -2.78%  "
$ws.Range("D15").Value = "0.523"
$ws.Range("E15").Value = "  -3.97%  "
$ws.Range("D16").Value = "25.986.55"
$ws.Range("E16").Value = "  -0.24%  "
$ws.Range("D17").Value = "0.0₃0740"
$ws.Range("E17").Value = "  -3.19%  "
$ws.Range("D18").Value = "61.24"
$ws.Range("E18").Value = "  -3.38%  "
$ws.Range("E19").Value = "  +0.42%  "
$ws.Range("D20").Value = "189.73"
$ws.Range("E20").Value = "  -3.00%  "
$ws.Range("D21").Value = "4.23"
$ws.Range("E21").Value = "  -2.86%  "
$ws.Range("D22").Value = "9.57"
$ws.Range("E22").Value = "  -3.65%  "
$ws.Range("D23").Value = "6.07"
$ws.Range("E23").Value = "  -2.59%  "
$ws.Range("D24").Value = "0.132"
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("D25").Value = "143.60"
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").Value = "1.78"
$ws.Range("E26").Value = "  -0.85%  "
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("D28").Value = "6.71"
$ws.Range("E28").Value = "  -2.78%  "
$ws.Range("D29").Value = "15.10"
$ws.Range("E29").Value = "  -2.90%  "
$ws.Range("E30").Value = "  -1.53%  "
$ws.Range("D31").Value = "0.0480"
$ws.Range("E31").Value = "  -3.52%  "
$ws.Range("D32").Value = "3.13"
$ws.Range("E32").Value = "  -4.31%  "
$ws.Range("E33").Value = "  -5.39%  "
$ws.Range("D34").Value = "2.41"
$ws.Range("E34").Value = "  -2.24%  "
$ws.Range("E35").Value = "  -3.71%  "
$ws.Range("D36").Value = "1.131.59"
$ws.Range("E36").Value = "  -0.24%  "
$ws.Range("D37").Value = "0.849"
$ws.Range("E37").Value = "  -6.43%  "
$ws.Range("E38").Value = "  -1.12%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.0154"
$ws.Range("E39").Value = "  -1.73%  "
$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").Value = "0.514"
$ws.Range("E40").Value = "  -5.04%  "
$ws.Range("D41").Value = "98.00"
$ws.Range("E41").Value = "  -1.34%  "
$ws.Range("E42").Value = "  -3.03%  "
$ws.Range("D43").Value = "1.769.57"
$ws.Range("E43").Value = "  -0.66%  "
$ws.Range("D44").Value = "5.19"
$ws.Range("E44").Value = "  -5.57%  "
$ws.Range("E45").Value = "  -2.50%  "
$ws.Range("D46").Value = "54.61"
$ws.Range("E46").Value = "  -3.89%  "
$ws.Range("E47").Value = "  +0.13%  "
$ws.Range("E48").Value = "  +0.58%  "
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("E50").Value = "  +0.56%  "
$ws.Range("E51").Value = "  -3.86%  "
